$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute")

# Insert a new row above row 3 (shifts rows 3+ down by one)
$ws.Rows("3:3").Insert()

# New row 3 gets the definition that used to live in row 2 (sample_id),
# plus a new attribute_name "sample_row_id"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 1).Value2 = "sample_row_id"
$ws.Cells.Item(3, 2).Value2 = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(3, 3).Value2 = "ordinal "
$ws.Cells.Item(3, 5).Value2 = "text"

# Row 2 (sample_id) gets the new, longer definition text
$ws.Cells.Item(2, 2).Value2 = "unique identifier for trap sample table based on julian date and year. Indicates if sample was included in analysis. Samples excluded from analysis contain a .0 "

# Row heights
$ws.Rows(2).RowHeight = 136
$ws.Rows(3).RowHeight = 51

# Update selection to match the author's final cursor position
$ws.Range("E2").Select()
